# leetcode 378 GFG array q6
# Adds two new rows (20, 21) to the "Leetcode" progress tracker for the
# "Array Order Statistics" GFG series (#1 and #6 / LeetCode 215 & 378),
# and updates the view's selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: GFG "Array Order Statistics #1" (LeetCode 215, kth largest) ---
$ws.Range("A20").Value = 215
$ws.Range("B20").Value = "Medium"
$ws.Range("C20").Value = "done"
$ws.Range("F20").Value = "Array Order Statistics #1"
$ws.Range("G20").Value = "https://www.geeksforgeeks.org/kth-smallestlargest-element-unsorted-array/"

# --- Row 21: GFG "Array Order Statistics #6" (LeetCode 378, kth smallest in sorted matrix) ---
$ws.Range("A21").Value = 378
$ws.Range("B21").Value = "Medium"
$ws.Range("C21").Value = "done"
$ws.Range("F21").Value = "Array Order Statistics #6"

# G21 carries a real hyperlink (like the other "Link" column entries), so add
# it via Hyperlinks.Add (no explicit display text -> Excel omits it, same as
# the other hyperlink rows) and then restyle it with the workbook's built-in
# Hyperlink cell style.
[void]$ws.Hyperlinks.Add($ws.Range("G21"), "https://www.geeksforgeeks.org/kth-smallest-element-in-a-row-wise-and-column-wise-sorted-2d-array-set-1/")
$ws.Range("G21").Style = "Hyperlink"

# --- View state: scroll the frozen pane down and move the active selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1
[void]$ws.Range("B22").Select()

Write-Output "Added rows 20-21 (Array Order Statistics #1 / #6) and updated selection."
